$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) cells are treated as text so values like
# trailing-zero decimals ("91.70", "1.010") are preserved verbatim,
# matching the source data which stores prices as text strings.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.566.94"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.923.77"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("E4").Value = "  +0.43%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.28"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4819"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4059"
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("E9").Value = "  +0.81%  "
$ws.Range("E10").Value = "  -0.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.88"
$ws.Range("E11").Value = "  +1.78%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.124"
$ws.Range("E12").Value = "  +1.79%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.902.16"
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.301"
$ws.Range("E14").Value = "  +2.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.70"
$ws.Range("E15").Value = "  +1.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06871"
$ws.Range("E16").Value = "  +1.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.013"
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001040"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.010"
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.574.13"
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.674"
$ws.Range("E22").Value = "  +0.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.02"
$ws.Range("E23").Value = "  +2.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.188"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.118.54"
$ws.Range("E25").Value = "  -0.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.03"
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.406"
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.090"
$ws.Range("E29").Value = "  -0.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.77"
$ws.Range("E30").Value = "  +0.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.012"
$ws.Range("E31").Value = "  -1.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09603"
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.616"
$ws.Range("E33").Value = "  +1.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.568"
$ws.Range("E34").Value = "  +0.15%  "
$ws.Range("E35").Value = "  -0.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06366"
$ws.Range("E36").Value = "  +4.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02290"
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.194"
$ws.Range("E38").Value = "  +1.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5954"
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("E40").Value = "  -0.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.011"
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.869"
$ws.Range("E42").Value = "  -1.78%  "
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("E44").Value = "  +1.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.286"
$ws.Range("E45").Value = "  +0.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.40"
$ws.Range("E46").Value = "  -0.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.07540"
$ws.Range("E47").Value = "  -0.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5556"
$ws.Range("E48").Value = "  -0.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.989"
$ws.Range("E49").Value = "  +2.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "119.28"
$ws.Range("E50").Value = "  +2.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.437"
$ws.Range("E51").Value = "  +0.78%  "
